$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header format from the existing last header cell (G1) onto the
# new header cell (H1) so it reuses the same cell style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New data value for the added column in the data row.
$ws.Range("H2").Value = 0
